{"js": "// Replace the 25 \"NN\u00f7N=Q, R\" answer strings in the practice-sheet table\n// with the new values from the commit, in document order. The very first\n// paragraph (the date heading, e.g. \"2024-08-19 Monday\") is left untouched;\n// every other non-empty paragraph in the document body corresponds, in\n// order, to one of the answer cells being updated.\n\nconst replacements = [\n  \"22\u00f77=3, 1\",\n  \"13\u00f73=4, 1\",\n  \"25\u00f79=2, 7\",\n  \"65\u00f76=10, 5\",\n  \"65\u00f72=32, 1\",\n  \"77\u00f77=11, 0\",\n  \"59\u00f78=7, 3\",\n  \"71\u00f75=14, 1\",\n  \"64\u00f78=8, 0\",\n  \"60\u00f76=10, 0\",\n  \"58\u00f78=7, 2\",\n  \"55\u00f77=7, 6\",\n  \"53\u00f72=26, 1\",\n  \"74\u00f78=9, 2\",\n  \"39\u00f79=4, 3\",\n  \"81\u00f78=10, 1\",\n  \"99\u00f74=24, 3\",\n  \"72\u00f78=9, 0\",\n  \"48\u00f72=24, 0\",\n  \"14\u00f78=1, 6\",\n  \"35\u00f76=5, 5\",\n  \"34\u00f73=11, 1\",\n  \"39\u00f72=19, 1\",\n  \"10\u00f73=3, 1\",\n  \"32\u00f78=4, 0\",\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Skip the heading paragraph (index 0) and any empty paragraphs (the blank\n// spacer rows in the table); apply the replacements, in order, to the\n// remaining paragraphs that hold an answer string.\nlet idx = 0;\nfor (let i = 1; i < paragraphs.items.length && idx < replacements.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text === \"\") continue;\n  para.insertText(replacements[idx], \"Replace\");\n  idx++;\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 \"NN\u00f7N=Q, R\" answer strings in the practice-sheet table\n# with the new values from the commit. The table is a 5-column grid where\n# only every 4th row (1, 5, 9, 13, 17) actually holds answer text - the\n# rows in between are blank spacer rows. The heading paragraph above the\n# table (\"2024-08-19 Monday\") is left untouched.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# New values, laid out exactly like the table: 5 data rows x 5 columns.\n$newValues = @(\n  @(\"22\u00f77=3, 1\", \"13\u00f73=4, 1\", \"25\u00f79=2, 7\", \"65\u00f76=10, 5\", \"65\u00f72=32, 1\"),\n  @(\"77\u00f77=11, 0\", \"59\u00f78=7, 3\", \"71\u00f75=14, 1\", \"64\u00f78=8, 0\", \"60\u00f76=10, 0\"),\n  @(\"58\u00f78=7, 2\", \"55\u00f77=7, 6\", \"53\u00f72=26, 1\", \"74\u00f78=9, 2\", \"39\u00f79=4, 3\"),\n  @(\"81\u00f78=10, 1\", \"99\u00f74=24, 3\", \"72\u00f78=9, 0\", \"48\u00f72=24, 0\", \"14\u00f78=1, 6\"),\n  @(\"35\u00f76=5, 5\", \"34\u00f73=11, 1\", \"39\u00f72=19, 1\", \"10\u00f73=3, 1\", \"32\u00f78=4, 0\")\n)\n\n$dataRows = @(1, 5, 9, 13, 17)\n\nfor ($i = 0; $i -lt $dataRows.Count; $i++) {\n  $r = $dataRows[$i]\n  for ($c = 1; $c -le 5; $c++) {\n    $cell = $t.Cell($r, $c)\n    $cell.Range.Text = $newValues[$i][$c - 1]\n  }\n}\n"}
